$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.1964449999999998
$ws.Range("I2").Value = 0.05113520435363903
$ws.Range("J2").Value = 0.05113520435363903
$ws.Range("M2").Value = 1.819857
$ws.Range("N2").Value = 5.459571
$ws.Range("O2").Value = 0.014853174625846068
$ws.Range("P2").Value = 0.014853174625846068
$ws.Range("Q2").Value = 1.332405269455
$ws.Range("R2").Value = 11.991647425095
$ws.Range("S2").Value = 0.0007595201197929247
$ws.Range("T2").Value = 0.0007595201197929247
$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.1964449999999998
$ws.Range("I3").Value = 0.05113520435363903
$ws.Range("J3").Value = 0.05113520435363903
$ws.Range("O3").Value = 0.7266185723345231
$ws.Range("P3").Value = 0.7266185723345231
$ws.Range("Q3").Value = 65.18137967473332
$ws.Range("R3").Value = 586.6324170725999
$ws.Range("S3").Value = 0.037155789183475285
$ws.Range("T3").Value = 0.037155789183475285
$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.1964449999999998
$ws.Range("I4").Value = 0.05113520435363903
$ws.Range("J4").Value = 0.05113520435363903
$ws.Range("M4").Value = 31.52924033333333
$ws.Range("N4").Value = 94.58772099999999
$ws.Range("O4").Value = 0.25733302808477204
$ws.Range("P4").Value = 0.25733302808477204
$ws.Range("Q4").Value = 23.08408076131611
$ws.Range("R4").Value = 207.75672685184495
$ws.Range("S4").Value = 0.01315877697805555
$ws.Range("T4").Value = 0.01315877697805555
$ws.Range("G5").Value = 0.7321483333333333
$ws.Range("H5").Value = 2.1964449999999998
$ws.Range("I5").Value = 0.05113520435363903
$ws.Range("J5").Value = 0.05113520435363903
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1464426666666667
$ws.Range("N5").Value = 0.43932800000000005
$ws.Range("O5").Value = 0.0011952249548588527
$ws.Range("P5").Value = 0.0011952249548588527
$ws.Range("Q5").Value = 0.1072177543288889
$ws.Range("R5").Value = 0.96495978896
$ws.Range("S5").Value = 0.00006111807231527643
$ws.Range("T5").Value = 0.00006111807231527643
$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("M6").Value = 1.819857
$ws.Range("N6").Value = 5.459571
$ws.Range("O6").Value = 0.014853174625846068
$ws.Range("P6").Value = 0.014853174625846068
$ws.Range("Q6").Value = 18.671440429641994
$ws.Range("R6").Value = 168.04296386677794
$ws.Range("S6").Value = 0.01064340932667491
$ws.Range("T6").Value = 0.01064340932667491
$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("O7").Value = 0.7266185723345231
$ws.Range("P7").Value = 0.7266185723345231
$ws.Range("S7").Value = 0.5206764940515162
$ws.Range("T7").Value = 0.5206764940515162
$ws.Range("I8").Value = 0.7165747117895102
$ws.Range("J8").Value = 0.7165747117895102
$ws.Range("M8").Value = 31.52924033333333
$ws.Range("N8").Value = 94.58772099999999
$ws.Range("O8").Value = 0.25733302808477204
$ws.Range("P8").Value = 0.25733302808477204
$ws.Range("Q8").Value = 323.484940122053
$ws.Range("R8").Value = 2911.3644610984766
$ws.Range("S8").Value = 0.18439834043376746
$ws.Range("T8").Value = 0.18439834043376746
$ws.Range("I9").Value = 0.7165747117895102
$ws.Range("J9").Value = 0.7165747117895102
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1464426666666667
$ws.Range("N9").Value = 0.43932800000000005
$ws.Range("O9").Value = 0.0011952249548588527
$ws.Range("P9").Value = 0.0011952249548588527
$ws.Range("Q9").Value = 1.5024782315448888
$ws.Range("R9").Value = 13.522304083903997
$ws.Range("S9").Value = 0.0008564679775516127
$ws.Range("T9").Value = 0.0008564679775516127
$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704000999999999
$ws.Range("I10").Value = 0.17935603462669877
$ws.Range("J10").Value = 0.17935603462669877
$ws.Range("M10").Value = 1.819857
$ws.Range("N10").Value = 5.459571
$ws.Range("O10").Value = 0.014853174625846068
$ws.Range("P10").Value = 0.014853174625846068
$ws.Range("Q10").Value = 4.6733933826189995
$ws.Range("R10").Value = 42.060540443571
$ws.Range("S10").Value = 0.002664006502509651
$ws.Range("T10").Value = 0.002664006502509651
$ws.Range("G11").Value = 2.568000333333333
$ws.Range("H11").Value = 7.704000999999999
$ws.Range("I11").Value = 0.17935603462669877
$ws.Range("J11").Value = 0.17935603462669877
$ws.Range("O11").Value = 0.7266185723345231
$ws.Range("P11").Value = 0.7266185723345231
$ws.Range("Q11").Value = 228.62280375585328
$ws.Range("R11").Value = 2057.60523380268
$ws.Range("S11").Value = 0.13032342582003315
$ws.Range("T11").Value = 0.13032342582003315
$ws.Range("G12").Value = 2.568000333333333
$ws.Range("H12").Value = 7.704000999999999
$ws.Range("I12").Value = 0.17935603462669877
$ws.Range("J12").Value = 0.17935603462669877
$ws.Range("M12").Value = 31.52924033333333
$ws.Range("N12").Value = 94.58772099999999
$ws.Range("O12").Value = 0.25733302808477204
$ws.Range("P12").Value = 0.25733302808477204
$ws.Range("Q12").Value = 80.96709968574676
$ws.Range("R12").Value = 728.7038971717208
$ws.Range("S12").Value = 0.04615423149576562
$ws.Range("T12").Value = 0.04615423149576562
$ws.Range("G13").Value = 2.568000333333333
$ws.Range("H13").Value = 7.704000999999999
$ws.Range("I13").Value = 0.17935603462669877
$ws.Range("J13").Value = 0.17935603462669877
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1464426666666667
$ws.Range("N13").Value = 0.43932800000000005
$ws.Range("O13").Value = 0.0011952249548588527
$ws.Range("P13").Value = 0.0011952249548588527
$ws.Range("Q13").Value = 0.3760648168142222
$ws.Range("R13").Value = 3.384583351328
$ws.Range("S13").Value = 0.00021437080839035888
$ws.Range("T13").Value = 0.00021437080839035888
$ws.Range("G14").Value = 0.7579039999999999
$ws.Range("H14").Value = 2.2737119999999997
$ws.Range("I14").Value = 0.05293404923015204
$ws.Range("J14").Value = 0.05293404923015204
$ws.Range("M14").Value = 1.819857
$ws.Range("N14").Value = 5.459571
$ws.Range("O14").Value = 0.014853174625846068
$ws.Range("P14").Value = 0.014853174625846068
$ws.Range("Q14").Value = 1.3792768997279998
$ws.Range("R14").Value = 12.413492097552
$ws.Range("S14").Value = 0.0007862386768685809
$ws.Range("T14").Value = 0.0007862386768685809
$ws.Range("G15").Value = 0.7579039999999999
$ws.Range("H15").Value = 2.2737119999999997
$ws.Range("I15").Value = 0.05293404923015204
$ws.Range("J15").Value = 0.05293404923015204
$ws.Range("O15").Value = 0.7266185723345231
$ws.Range("P15").Value = 0.7266185723345231
$ws.Range("Q15").Value = 67.47434383423999
$ws.Range("R15").Value = 607.2690945081599
$ws.Range("S15").Value = 0.03846286327949844
$ws.Range("T15").Value = 0.03846286327949844
$ws.Range("G16").Value = 0.7579039999999999
$ws.Range("H16").Value = 2.2737119999999997
$ws.Range("I16").Value = 0.05293404923015204
$ws.Range("J16").Value = 0.05293404923015204
$ws.Range("M16").Value = 31.52924033333333
$ws.Range("N16").Value = 94.58772099999999
$ws.Range("O16").Value = 0.25733302808477204
$ws.Range("P16").Value = 0.25733302808477204
$ws.Range("Q16").Value = 23.89613736559466
$ws.Range("R16").Value = 215.06523629035195
$ws.Range("S16").Value = 0.01362167917718342
$ws.Range("T16").Value = 0.01362167917718342
$ws.Range("G17").Value = 0.7579039999999999
$ws.Range("H17").Value = 2.2737119999999997
$ws.Range("I17").Value = 0.05293404923015204
$ws.Range("J17").Value = 0.05293404923015204
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1464426666666667
$ws.Range("N17").Value = 0.43932800000000005
$ws.Range("O17").Value = 0.0011952249548588527
$ws.Range("P17").Value = 0.0011952249548588527
$ws.Range("Q17").Value = 0.11098948283733334
$ws.Range("R17").Value = 0.998905345536
$ws.Range("S17").Value = 0.00006326809660160476
$ws.Range("T17").Value = 0.00006326809660160476
